# Penalty/Reward System attempt (unfinished):
# - On "Weekly Quantity" sheet, the weekly bucket dated 2024-03-10 (serial 45368.99999999999,
#   quantity 60) is removed entirely, shifting all later weeks up by one row.
# - On "Monthly Trend" sheet, the March 2024 monthly total (row 9) is reduced from 820 to 760
#   to reflect the removed 60-unit weekly entry.

$wb = $excel.ActiveWorkbook

$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")

# Remove the row for week 45368.99999999999 (qty 60); remaining rows shift up automatically.
$wsWeekly.Range("A23").EntireRow.Delete()

# Update the corresponding monthly total.
$wsMonthly.Range("B9").Value = 760
